$p = $ppt.ActivePresentation
$s = $p.Slides.Item(31)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph describing PUSH order ("Always pushes registers in same order") ---
$pushPara = $tr.Paragraphs(6, 1)
# Reset through a "neutral" placeholder first so the COM text-setter's
# common prefix/suffix preservation can't accidentally keep a fragment
# of the old run boundary-adjacent to the new text.
$pushPara.Text = "0123456789"
$pushText = "Largest register number pushed first (to largest address)"
$pushPara.Text = $pushText

# Bold the leading "Largest "
$pushPara.Characters(1, 8).Font.Bold = $true
# Bold the second occurrence of "largest "
$idx = $pushText.IndexOf("largest ")
$pushPara.Characters($idx + 1, 8).Font.Bold = $true

# --- Paragraph describing POP order ("Always pops registers in same order ...") ---
$popPara = $tr.Paragraphs(13, 1)
$popPara.Text = "0123456789"
$popText = "Smallest register number popped first (from smallest address)"
$popPara.Text = $popText

# Bold the leading "Smallest "
$popPara.Characters(1, 9).Font.Bold = $true
# Bold the second occurrence of "smallest "
$idx2 = $popText.IndexOf("smallest ")
$popPara.Characters($idx2 + 1, 9).Font.Bold = $true
